# Add a new row of contributor details (name / email / repo link) to the
# "Open Source task" worksheet, wire up a mailto hyperlink on the email
# cell (matching the style already used by the existing email/link row),
# and leave the active selection on C5 - mirroring the author's commit
# "Added my details to the Excel file".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 3): name, email, repo link.
$ws.Range("A3").Value = "محمد حسين غنيم طوخي"
$ws.Range("B3").Value = "eng.nooone@gmail.com"
$ws.Range("C3").Value = "https://github.com/EngNoOne/Security-Task.git"

# Hyperlink the email cell, same as the existing B2 email hyperlink.
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:eng.nooone@gmail.com")

# Hyperlinks.Add() mints its own cell style - reapply the workbook's
# built-in "Hyperlink" style so B3 reuses the same style already applied
# to B2/C2 instead of drifting from it.
$ws.Range("B3").Style = "Hyperlink"

# Leave the selection where the author left it after entering the data.
$ws.Range("C5").Select()
